$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.896.03"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "1.631.31"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.37"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.49"
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("E10").Value = "  -0.23%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "1.862.70"
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").Value = "1.621.92"
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("E14").Value = "  -1.34%  "
$ws.Range("E15").Value = "  -1.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.36"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").Value = "27.899.25"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.04"
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("E23").Value = "  -3.54%  "
$ws.Range("E24").Value = "  -0.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.72"
$ws.Range("E25").Value = "  +1.10%  "
$ws.Range("E26").Value = "  -0.80%  "
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0482"
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("E33").Value = "  +0.57%  "
$ws.Range("D34").Value = "1.394.61"
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("E36").Value = "  +10.20%  "
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("E38").Value = "  +1.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.559"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.847"
$ws.Range("E40").Value = "  -3.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  -0.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.83"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "65.87"
$ws.Range("E44").Value = "  -1.88%  "
$ws.Range("E45").Value = "  -1.82%  "
$ws.Range("D46").Value = "1.772.91"
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.15"
$ws.Range("E47").Value = "  -3.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.65"
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("E49").Value = "  +1.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0505"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.65"
$ws.Range("E51").Value = "  +0.60%  "
